$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -5
    13 = -7
    17 = -3
    23 = 0
    24 = 2
    28 = 5
    29 = 4
    32 = -3
    35 = -3
    43 = 4
    44 = 5
    46 = -3
    56 = -2
    58 = -8
    64 = 11
    67 = 1
    68 = -7
    70 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
